$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 36, column A: was stored as text "71277620" -> change to numeric 71277620
$ws.Cells.Item(36, 1).Value = 71277620

# New row 37: add payment 71277620 (Cash) 2025-08-18T17:11:19
$ws.Cells.Item(37, 1).Value = "'71277620"
$ws.Cells.Item(37, 2).Value = ""
$ws.Cells.Item(37, 3).Value = "Cash"
$ws.Cells.Item(37, 4).Value = "2025-08-18T17:11:19"
$ws.Cells.Item(37, 5).Value = 76
$ws.Cells.Item(37, 6).Value = ""
$ws.Cells.Item(37, 7).Value = 76
$ws.Cells.Item(37, 8).Value = 0
$ws.Cells.Item(37, 9).Value = 0
